$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 first (higher row number), then row 26, so row indices
# for the earlier deletion stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
